# Update the BOM worksheet:
#  - Row 7 (Q3): part changes from IRLM2502/C347503 to IRLM6402/C347504
#  - Old row 8 (J3, Conn_02x20_Odd_Even header) is removed entirely,
#    shifting the following rows (Q2, J1/Micro USB) up by one
#  - Selection moves to C14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Comment/Part# for the Q3 row (row 7)
$ws.Range("A7").Value = "IRLM6402"
$ws.Range("D7").Value = "C347504"

# Remove the obsolete connector row (old row 8: J3 / Conn_02x20_Odd_Even)
$ws.Rows(8).Delete() | Out-Null

# Update the active selection to match the saved view state
$ws.Range("C14").Select() | Out-Null
